$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on column D data rows to prevent Excel from
# auto-converting numeric-looking strings (e.g. "7.51") into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.896.91"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.632.12"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "23.21"
$ws.Range("E8").Value = "  -1.16%  "
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "1.864.03"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "1.634.52"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "65.20"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "27.896.34"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "229.76"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "10.35"
$ws.Range("E23").Value = "  -3.03%  "
$ws.Range("E24").Value = "  -4.27%  "
$ws.Range("D25").Value = "153.96"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "15.61"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("D34").Value = "1.395.78"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +9.47%  "
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "0.560"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "66.84"
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("D44").Value = "5.53"
$ws.Range("E44").Value = "  +3.08%  "
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("E46").Value = "  -1.19%  "
$ws.Range("D47").Value = "1.773.28"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "87.63"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").Value = "  -0.24%  "

# Restore default (Normal) style on column D so no stray number format
# is left behind on cells that did not need one.
$ws.Range("D2:D51").Style = "Normal"
